$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I, row 2: empty cell, same formatting (border/font) as H2
$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial(-4122)

# Column I, row 3: header year 2021, same formatting as H3 but font size 11
$ws.Range("H3").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("I3").Value = 2021
$ws.Range("I3").Font.Size = 11

# Column I, row 4: value 149, same formatting as H4 but font size 11
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = 149
$ws.Range("I4").Font.Size = 11

# Column I, row 5: value 159, same formatting as H5 but font size 11
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").Value = 159
$ws.Range("I5").Font.Size = 11

# Clear clipboard marching ants / move selection like the source file
$ws.Range("K4").Select()
